$d = $word.ActiveDocument

# 1. "döpte till Menu()" -> "döpte till MainMenu()"
$d.Content.Find.Execute("döpte till Menu()", $true, $false, $false, $false, $false,
                         $true, 1, $false, "döpte till MainMenu()", 2)

# 2. "Jag skapade metoden Read()" -> "Jag skapade metoden ReadFromTextFile()"
$d.Content.Find.Execute("Jag skapade metoden Read()", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Jag skapade metoden ReadFromTextFile()", 2)

# 2b. "metoden AddMember()" -> "metoden AddMemberToList()"
$d.Content.Find.Execute("metoden AddMember()", $true, $false, $false, $false, $false,
                         $true, 1, $false, "metoden AddMemberToList()", 2)

# 2c. "Read() anropas" -> "ReadFromTextFile() anropas"
$d.Content.Find.Execute("Read() anropas", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ReadFromTextFile() anropas", 2)

# 3. "metod, Details()" -> "metod, MemberDetails()"
$d.Content.Find.Execute("metod, Details()", $true, $false, $false, $false, $false,
                         $true, 1, $false, "metod, MemberDetails()", 2)

# 4. "metoden Describe()" -> "metoden DescribeMember()"
$d.Content.Find.Execute("metoden Describe()", $true, $false, $false, $false, $false,
                         $true, 1, $false, "metoden DescribeMember()", 2)

# 5. "metoden Remove()" -> "metoden DeleteMember()"
$d.Content.Find.Execute("metoden Remove()", $true, $false, $false, $false, $false,
                         $true, 1, $false, "metoden DeleteMember()", 2)

# 6. "metoden Create()" -> "metoden CreateNewMember()"
$d.Content.Find.Execute("metoden Create()", $true, $false, $false, $false, $false,
                         $true, 1, $false, "metoden CreateNewMember()", 2)
